{"js": "// Bug fix: an empty footnote (body text is just a stray \"\u0964\" punctuation\n// mark with no real content) was left in the document together with its\n// reference mark; also a stray \"a\" had been appended to the end of a\n// different, legitimate footnote. This removes the empty note (and its\n// reference in the body) and strips the stray trailing \"a\".\n\nconst footnotes = context.document.body.footnotes;\nfootnotes.load(\"items\");\nawait context.sync();\n\n// Load every footnote body's text so we can find the two notes we need by\n// their content instead of relying on a brittle, hard-coded index.\nconst bodies = footnotes.items.map((fn) => fn.body);\nbodies.forEach((b) => b.load(\"text\"));\nawait context.sync();\n\nlet strayALetterIndex = -1;\nlet emptyNoteIndex = -1;\n\nfor (let i = 0; i < footnotes.items.length; i++) {\n  // Footnote body text starts with the footnote-mark control char (\\u0002)\n  // followed by a space and then the actual note text.\n  const noteText = bodies[i].text.replace(/^[\\s\\u0002]+/, \"\");\n\n  if (noteText.endsWith(\"\u0f54\u0f7a\u0f0b\u0f45\u0f72\u0f53\u0f0da\")) {\n    strayALetterIndex = i;\n  }\n  if (noteText.trim() === \"\u0f0d\") {\n    emptyNoteIndex = i;\n  }\n}\n\n// 1. Strip the stray trailing \"a\" from the legitimate footnote's text.\nif (strayALetterIndex !== -1) {\n  const target = footnotes.items[strayALetterIndex].body.search(\"\u0f54\u0f7a\u0f0b\u0f45\u0f72\u0f53\u0f0da\", {\n    matchCase: true,\n  });\n  target.load(\"items\");\n  await context.sync();\n\n  target.items[0].insertText(\"\u0f54\u0f7a\u0f0b\u0f45\u0f72\u0f53\u0f0d\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 2. Delete the empty footnote entirely, which also removes its reference\n// mark from the body paragraph.\nif (emptyNoteIndex !== -1) {\n  footnotes.items[emptyNoteIndex].reference.delete();\n  await context.sync();\n}\n", "ps1": "# Bug fix: an empty footnote (body text is just a stray \"\u0964\" punctuation\n# mark with no real content) was left in the document together with its\n# reference mark; also a stray \"a\" had been appended to the end of a\n# different, legitimate footnote. This removes the empty note (and its\n# reference in the body) and strips the stray trailing \"a\".\n\n$d = $word.ActiveDocument\n\n# Find the two notes we need by their content instead of relying on a\n# brittle, hard-coded index.\n$strayLetterIndex = -1\n$emptyNoteIndex = -1\n\nfor ($i = 1; $i -le $d.Footnotes.Count; $i++) {\n    $noteText = $d.Footnotes.Item($i).Range.Text.Trim()\n\n    if ($noteText.EndsWith(\"a\")) {\n        $strayLetterIndex = $i\n    }\n    if ($noteText -eq \"\u0f0d\") {\n        $emptyNoteIndex = $i\n    }\n}\n\n# 1. Strip the stray trailing \"a\" from the legitimate footnote's text.\nif ($strayLetterIndex -ne -1) {\n    $rng = $d.Footnotes.Item($strayLetterIndex).Range\n    $rng.Text = $rng.Text.Substring(0, $rng.Text.Length - 1)\n}\n\n# 2. Delete the empty footnote entirely, which also removes its reference\n# mark from the body paragraph.\nif ($emptyNoteIndex -ne -1) {\n    $d.Footnotes.Item($emptyNoteIndex).Delete()\n}\n"}
